$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix Razon social / Nombre Fantasia punctuation (commas -> periods) ---
$ws.Range("E84").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E159").Value = "PARPAGNOLI. PEDRO RICARDO"
$ws.Range("F159").Value = "PARPAGNOLI. PEDRO RICARDO"
$ws.Range("E181").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"
$ws.Range("E194").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"

# --- Fix Importe column: "1.234,56" (es-AR) -> "1234.56" (text, dot decimal) ---
$importe = @(
    @{Row=2; Value="19020.00"},
    @{Row=3; Value="22500.00"},
    @{Row=4; Value="11730.00"},
    @{Row=5; Value="11000.00"},
    @{Row=6; Value="2779.00"},
    @{Row=7; Value="9420.00"},
    @{Row=8; Value="160000.00"},
    @{Row=9; Value="4390.00"},
    @{Row=10; Value="54716.20"},
    @{Row=11; Value="430760.00"},
    @{Row=12; Value="2150.00"},
    @{Row=13; Value="450.00"},
    @{Row=14; Value="812859.40"},
    @{Row=15; Value="23161.60"},
    @{Row=16; Value="10500.00"},
    @{Row=17; Value="267376.61"},
    @{Row=18; Value="457500.00"},
    @{Row=19; Value="339105.80"},
    @{Row=20; Value="39075.00"},
    @{Row=21; Value="93680.00"},
    @{Row=22; Value="33867.20"},
    @{Row=23; Value="1229.03"},
    @{Row=24; Value="9975.00"},
    @{Row=25; Value="1365.00"},
    @{Row=26; Value="2490.00"},
    @{Row=27; Value="11200.00"},
    @{Row=28; Value="768.00"},
    @{Row=29; Value="9440.00"},
    @{Row=30; Value="5902.96"},
    @{Row=31; Value="74.00"},
    @{Row=32; Value="1640.00"},
    @{Row=33; Value="109310.00"},
    @{Row=34; Value="507.36"},
    @{Row=35; Value="18270.10"},
    @{Row=36; Value="2009947.80"},
    @{Row=37; Value="26777.94"},
    @{Row=38; Value="220.00"},
    @{Row=39; Value="6573.99"},
    @{Row=40; Value="8550.00"},
    @{Row=41; Value="1477.00"},
    @{Row=42; Value="24474.00"},
    @{Row=43; Value="18755.00"},
    @{Row=44; Value="10400.00"},
    @{Row=45; Value="44209.73"},
    @{Row=46; Value="5200.00"},
    @{Row=47; Value="10015.00"},
    @{Row=48; Value="1834.80"},
    @{Row=49; Value="9716.79"},
    @{Row=50; Value="10800.00"},
    @{Row=51; Value="5300.00"},
    @{Row=52; Value="5881.00"},
    @{Row=53; Value="9940.00"},
    @{Row=54; Value="14945.62"},
    @{Row=55; Value="282000.00"},
    @{Row=56; Value="5000.00"},
    @{Row=57; Value="647.35"},
    @{Row=58; Value="1990.00"},
    @{Row=59; Value="7812.90"},
    @{Row=60; Value="3001.00"},
    @{Row=61; Value="2250.00"},
    @{Row=62; Value="43324.59"},
    @{Row=63; Value="1760.00"},
    @{Row=64; Value="28000.00"},
    @{Row=65; Value="59500.00"},
    @{Row=66; Value="29910.00"},
    @{Row=67; Value="65600.00"},
    @{Row=68; Value="7200.00"},
    @{Row=69; Value="9735.90"},
    @{Row=70; Value="3612.00"},
    @{Row=71; Value="144000.00"},
    @{Row=72; Value="1329.99"},
    @{Row=73; Value="8400.00"},
    @{Row=74; Value="5922.40"},
    @{Row=75; Value="18000.00"},
    @{Row=76; Value="3100.00"},
    @{Row=77; Value="3260.00"},
    @{Row=78; Value="4700.00"},
    @{Row=79; Value="780.00"},
    @{Row=80; Value="15500.00"},
    @{Row=81; Value="228750.00"},
    @{Row=82; Value="10836.00"},
    @{Row=83; Value="9700.00"},
    @{Row=84; Value="1600.00"},
    @{Row=85; Value="8280.00"},
    @{Row=86; Value="2438.00"},
    @{Row=87; Value="517.29"},
    @{Row=88; Value="814.90"},
    @{Row=89; Value="6330.00"},
    @{Row=90; Value="4194.00"},
    @{Row=91; Value="9.39"},
    @{Row=92; Value="5.95"},
    @{Row=93; Value="4.00"},
    @{Row=94; Value="4.52"},
    @{Row=95; Value="8064.21"},
    @{Row=96; Value="22266.11"},
    @{Row=97; Value="210.00"},
    @{Row=98; Value="1179.75"},
    @{Row=99; Value="15120.00"},
    @{Row=100; Value="4162.95"},
    @{Row=101; Value="2934.00"},
    @{Row=102; Value="4117.52"},
    @{Row=103; Value="31983.00"},
    @{Row=104; Value="25240.00"},
    @{Row=105; Value="505.00"},
    @{Row=106; Value="11730.00"},
    @{Row=107; Value="10140.00"},
    @{Row=108; Value="11510.00"},
    @{Row=109; Value="22888.97"},
    @{Row=110; Value="3469.98"},
    @{Row=111; Value="14485.01"},
    @{Row=112; Value="865.80"},
    @{Row=113; Value="345.00"},
    @{Row=114; Value="53183.00"},
    @{Row=115; Value="800.00"},
    @{Row=116; Value="2294.70"},
    @{Row=117; Value="6300.00"},
    @{Row=118; Value="768000.00"},
    @{Row=119; Value="937.00"},
    @{Row=120; Value="6360.00"},
    @{Row=121; Value="25400.00"},
    @{Row=122; Value="33000.00"},
    @{Row=123; Value="7000.00"},
    @{Row=124; Value="4000.00"},
    @{Row=125; Value="4200.00"},
    @{Row=126; Value="25000.00"},
    @{Row=127; Value="8000.00"},
    @{Row=128; Value="23000.00"},
    @{Row=129; Value="185000.00"},
    @{Row=130; Value="2173.09"},
    @{Row=131; Value="3651.10"},
    @{Row=132; Value="23372.00"},
    @{Row=133; Value="25430.00"},
    @{Row=134; Value="6061.00"},
    @{Row=135; Value="1560.00"},
    @{Row=136; Value="6510000.00"},
    @{Row=137; Value="50925.00"},
    @{Row=138; Value="25000.00"},
    @{Row=139; Value="16000.00"},
    @{Row=140; Value="36000.00"},
    @{Row=141; Value="8000.00"},
    @{Row=142; Value="13000.00"},
    @{Row=143; Value="5000.00"},
    @{Row=144; Value="17017.00"},
    @{Row=145; Value="12000.00"},
    @{Row=146; Value="14000.00"},
    @{Row=147; Value="10000.00"},
    @{Row=148; Value="6000.00"},
    @{Row=149; Value="12000.00"},
    @{Row=150; Value="10000.00"},
    @{Row=151; Value="10000.00"},
    @{Row=152; Value="16000.00"},
    @{Row=153; Value="12000.00"},
    @{Row=154; Value="14000.00"},
    @{Row=155; Value="12500.00"},
    @{Row=156; Value="12000.00"},
    @{Row=157; Value="16000.00"},
    @{Row=158; Value="3000.00"},
    @{Row=159; Value="6000.00"},
    @{Row=160; Value="16500.00"},
    @{Row=161; Value="6000.00"},
    @{Row=162; Value="118111.44"},
    @{Row=163; Value="15000.00"},
    @{Row=164; Value="13000.00"},
    @{Row=165; Value="12000.00"},
    @{Row=166; Value="14000.00"},
    @{Row=167; Value="61153.40"},
    @{Row=168; Value="7000.00"},
    @{Row=169; Value="50000.00"},
    @{Row=170; Value="6000.00"},
    @{Row=171; Value="8600.00"},
    @{Row=172; Value="8800.00"},
    @{Row=173; Value="3700.00"},
    @{Row=174; Value="11020.00"},
    @{Row=175; Value="11800.00"},
    @{Row=176; Value="617412.04"},
    @{Row=177; Value="9864.00"},
    @{Row=178; Value="10400.00"},
    @{Row=179; Value="25300.76"},
    @{Row=180; Value="1671.40"},
    @{Row=181; Value="3055.00"},
    @{Row=182; Value="1820.50"},
    @{Row=183; Value="6800.00"},
    @{Row=184; Value="4200.00"},
    @{Row=185; Value="19710.00"},
    @{Row=186; Value="2620.00"},
    @{Row=187; Value="5340.00"},
    @{Row=188; Value="1836.00"},
    @{Row=189; Value="810.00"},
    @{Row=190; Value="1600.00"},
    @{Row=191; Value="4300.00"},
    @{Row=192; Value="26594.33"},
    @{Row=193; Value="3218.60"},
    @{Row=194; Value="4970.00"},
    @{Row=195; Value="560.00"},
    @{Row=196; Value="2700.00"},
    @{Row=197; Value="9900.00"},
    @{Row=198; Value="9945.00"},
    @{Row=199; Value="1100.00"},
    @{Row=200; Value="15300.00"},
    @{Row=201; Value="9124.11"},
    @{Row=202; Value="480.00"},
    @{Row=203; Value="7023.95"},
    @{Row=204; Value="480000.00"},
    @{Row=205; Value="4730174.63"},
    @{Row=206; Value="2900.00"},
    @{Row=207; Value="4000.00"},
    @{Row=208; Value="4100.00"},
    @{Row=209; Value="6800.00"},
    @{Row=210; Value="37300.00"},
    @{Row=211; Value="11500.00"},
    @{Row=212; Value="3450.00"},
    @{Row=213; Value="57200.00"},
    @{Row=214; Value="60000.00"},
    @{Row=215; Value="21000.00"},
    @{Row=216; Value="9770.00"},
    @{Row=217; Value="525000.00"},
    @{Row=218; Value="65800.00"},
    @{Row=219; Value="2240000.00"},
    @{Row=220; Value="462874.62"},
    @{Row=221; Value="799012.50"},
    @{Row=222; Value="13050.00"},
    @{Row=223; Value="8550.00"},
    @{Row=224; Value="39000.00"},
    @{Row=225; Value="2872.86"},
    @{Row=226; Value="8000.00"},
    @{Row=227; Value="15600.00"},
    @{Row=228; Value="2000.00"},
    @{Row=229; Value="23220.00"}
)

foreach ($item in $importe) {
    $cell = $ws.Cells.Item($item.Row, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
    $cell.Style = "Normal"
}